$d = $word.ActiveDocument

# Find the index (1-based) of the "Docente(s) Responsável(eis) " heading paragraph
$idx = 1
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "Docente(s) Responsável(eis) ") {
        $targetIdx = $idx
        break
    }
    $idx = $idx + 1
}

# Insert a new empty paragraph right after it
$target = $d.Paragraphs.Item($targetIdx)
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# Fill in the new paragraph's text and style
$newPara = $d.Paragraphs.Item($targetIdx + 1)
$newPara.Range.Text = "6270264 - Juan Fernando Zapata Zapata"
$newPara.Style = "ListBullet"
